# Daily attendance processing - swap the order of "Recorded By" entries
# that list both the submitter e-mail and "System" in the
# "dnasr281@gmail.com, System" form, re-ordering them to "System, dnasr281@gmail.com".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

[void]$ws.Cells.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com")
